$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: "  Garibaldi meets King Vittorio Emanuele ..." ->
#           "  Giuseppe Garibaldi meets King Vittorio Emanuele ..."
# The original run "  Garibaldi meets King ..." must become three runs:
#   "  "  |  "Giuseppe "  |  "Garibaldi meets King ..."
# all sharing the same (color-only) run formatting.
# ---------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("Garibaldi meets King", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not find 'Garibaldi meets King' text"
}

$insertionPoint1 = $d.Range($rng1.Start, $rng1.Start)
$insertionPoint1.InsertBefore("Giuseppe ")

# Touch the formatting of the newly-inserted run (toggle Bold off/on and
# back) so it is emitted as its own <w:r> instead of being silently
# re-absorbed into a neighbouring run of identical formatting.
$newRun1 = $d.Range($rng1.Start, $rng1.Start + 9)
$newRun1.Bold = 1
$newRun1.Bold = 0

# ---------------------------------------------------------------------
# Change 2: "March 2016" -> "July 2016", split across two runs:
#   "July"  |  " 2016"
# ---------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("March 2016", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find 'March 2016' text"
}

$marchRng = $d.Range($rng2.Start, $rng2.Start + 5)
$marchRng.Text = "July"

$newRun2 = $d.Range($rng2.Start, $rng2.Start + 4)
$newRun2.Bold = 1
$newRun2.Bold = 0
